$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Socks in the Dark:" title -> bold + underline (paragraph mark and run)
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Socks in the Dark:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$para.Range.Font.Bold = 1
$para.Range.Font.Underline = 1

# ---------------------------------------------------------------------------
# 2) "Predicting Fingers:" -> bold + underline on "Predicting Fingers" only,
#    the trailing colon stays unformatted (separate run).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Predicting Fingers", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Font.Bold = 1
$rng.Font.Underline = 1

# ---------------------------------------------------------------------------
# 3) Move the "_GoBack" bookmark from after " solution" (Socks paragraph) to
#    inside "girl counting" ("girl co|unting") in the Predicting Fingers
#    description paragraph. Adding a bookmark with the same name relocates
#    it (Word only keeps one bookmark per name).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("girl co", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insPoint = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $insPoint)

# ---------------------------------------------------------------------------
# 4) Split the "...a minimum of  12 because..." run and wrap "of  12" with
#    proofErr gramStart/gramEnd markers, matching Word's grammar-check
#    annotations. Rebuild the whole paragraph via InsertXML so the proofErr
#    elements land exactly between the split runs.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -match "The first of the pair of socks") {
        $target = $p
        break
    }
}

$newParaXml = '<w:p w14:paraId="41A6F975" w14:textId="77777777" w:rsidR="00E13894" w:rsidRDefault="00E13894" w:rsidP="00E13894"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">The first of the pair of socks would take 4 socks because you could pull one of each sock and then pull one more that would complete the pair. The 3 pairs of all the colors would require a minimum </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>of  12</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">because if you take the odds of the first question to get 1 pair it would take the same amount 3 times which adds to 12 socks. </w:t></w:r></w:p>'

$pkg = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
  '<w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.Range.InsertXML($pkg)

Write-Host "Edits applied"
